# feat: add 2022-Q1 data
#
# - insert a new "2022-Q1" sheet (holdings detail) between "2021-Q4" and "总计"
# - update the "总计" (totals) summary sheet with a new leading row for 2022-Q1
#   and shift the existing 2021-Q4 / 2021-Q3 rows down

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Re-create the "总计" sheet in the right tab order so the new sheet lands
#    between "2021-Q4" and "总计", and sheetIds line up (2022-Q1 -> 3, 总计 -> 4)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)   # 2021-Q3
$ws2 = $wb.Worksheets.Item(2)   # 2021-Q4
$oldTotal = $wb.Worksheets.Item("总计")
$oldTotal.Delete()

$wsQ1 = $wb.Worksheets.Add($null, $ws2)
$wsQ1.Name = "2022-Q1"
$wsQ1.PageSetup.LeftMargin = 54
$wsQ1.PageSetup.RightMargin = 54
$wsQ1.PageSetup.TopMargin = 72
$wsQ1.PageSetup.BottomMargin = 72
$wsQ1.PageSetup.HeaderMargin = 36
$wsQ1.PageSetup.FooterMargin = 36

$wsTotal = $wb.Worksheets.Add($null, $wsQ1)
$wsTotal.Name = "总计"
$wsTotal.PageSetup.LeftMargin = 54
$wsTotal.PageSetup.RightMargin = 54
$wsTotal.PageSetup.TopMargin = 72
$wsTotal.PageSetup.BottomMargin = 72
$wsTotal.PageSetup.HeaderMargin = 36
$wsTotal.PageSetup.FooterMargin = 36

# ---------------------------------------------------------------------------
# 2. Populate "2022-Q1" - same layout as the other quarterly sheets.
#    Copy the header row + A2 index cell from "2021-Q4" first so the
#    bordered/bold/centered "s=2" formatting matches exactly, then overwrite
#    with the real 2022-Q1 values.
# ---------------------------------------------------------------------------
$ws2.Range("B1:H2").Copy($wsQ1.Range("B1"))
$ws2.Range("A2").Copy($wsQ1.Range("A2"))

$wsQ1.Range("B1").Value = "基金代码"
$wsQ1.Range("C1").Value = "基金名称"
$wsQ1.Range("D1").Value = "基金规模"
$wsQ1.Range("E1").Value = "股票总仓位"
$wsQ1.Range("F1").Value = "仓位占比"
$wsQ1.Range("G1").Value = "持有市值(亿元)"
$wsQ1.Range("H1").Value = "仓位排名"

# Text-valued columns B2:G2 must stay text (not auto-coerced to numbers) -
# apply a text format before writing, then clear the format again so the
# cells end up with the default (no explicit) style, matching the source data.
$q1TextRng = $wsQ1.Range("B2:G2")
$q1TextRng.NumberFormat = "@"
$wsQ1.Range("B2").Value = "001075"
$wsQ1.Range("C2").Value = "宝盈转型动力灵活配置混合"
$wsQ1.Range("D2").Value = "5.13"
$wsQ1.Range("E2").Value = "86.64"
$wsQ1.Range("F2").Value = "3.40"
$wsQ1.Range("G2").Value = "0.1744"
$q1TextRng.ClearFormats()

$wsQ1.Range("H2").Value = 6

# ---------------------------------------------------------------------------
# 3. Populate "总计" - header + 3 rows (2022-Q1, 2021-Q4, 2021-Q3).
#    Copy the header + index-column styling from "2021-Q4" (also "s=2").
# ---------------------------------------------------------------------------
$ws2.Range("B1:D1").Copy($wsTotal.Range("B1"))
$ws2.Range("A2").Copy($wsTotal.Range("A2"))
$ws2.Range("A2").Copy($wsTotal.Range("A3"))
$ws2.Range("A2").Copy($wsTotal.Range("A4"))

$wsTotal.Range("B1").Value = "日期"
$wsTotal.Range("C1").Value = "持有数量(只)"
$wsTotal.Range("D1").Value = "持有市值(亿元)"

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q1"
$wsTotal.Range("C2").Value = 1
$wsTotal.Range("D2").Value = 0.17

$wsTotal.Range("A3").Value = 1
$wsTotal.Range("B3").Value = "2021-Q4"
$wsTotal.Range("C3").Value = 1
$wsTotal.Range("D3").Value = 0.21

$wsTotal.Range("A4").Value = 2
$wsTotal.Range("B4").Value = "2021-Q3"
$wsTotal.Range("C4").Value = 1
$wsTotal.Range("D4").Value = 0.02

# Restore the originally-active tab/selection (unchanged by this edit).
$ws1.Activate()
$ws1.Range("A1").Select()
